$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")
$ws.Activate()

# Mark all remaining "Implemented"/"Failed" Testing-column results as "Done"
# (E74 was "Failed"; E76:E110, skipping the category header rows, were "Implemented")
$testingRows = @(74,76,77,78,80,81,82,83,84,85,86,88,89,90,91,93,94,95,96,97,98,99,100,101,102,103,104,106,107,108,109,110)
foreach ($r in $testingRows) {
    $ws.Range("E" + $r).Value = "Done"
}

# Reflect the saved view state: zoom level and current selection
$excel.ActiveWindow.Zoom = 125
$ws.Range("M6").Select() | Out-Null
